# PaymentsELHardCoded.xlsx - RAD Updates for August 30, Phase 2
# The "Date" column (B) on Sheet1 holds the timestamp recorded the last
# time each hard-coded payments test row was executed/generated by the
# Katalon RAD tooling. Refresh the four rows that were re-run, stamping
# them with their new execution timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value  = "Mon Jul 17 21:27:24 EDT 2023"
$ws.Range("B10").Value = "Mon Jul 17 21:28:09 EDT 2023"
$ws.Range("B11").Value = "Mon Jul 17 21:28:51 EDT 2023"
$ws.Range("B12").Value = "Mon Jul 17 21:29:32 EDT 2023"
